# Scheduled market-data refresh: update computed Leve profit columns (H:N) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 343577.1
$ws.Range("I15").Value = 343577.1
$ws.Range("K15").Value = 1030731.3
$ws.Range("M15").Value = -1030562.3

$ws.Range("H19").Value = 1261.1666
$ws.Range("I19").Value = 600
$ws.Range("J19").Value = 1343.8125
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 1343.8125
$ws.Range("M19").Value = -425
$ws.Range("N19").Value = -1693.8125

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H28").Value = 515.1905
$ws.Range("I28").Value = 485.4375
$ws.Range("J28").Value = 610.4
$ws.Range("K28").Value = 485.4375
$ws.Range("L28").Value = 610.4
$ws.Range("M28").Value = -0.4375
$ws.Range("N28").Value = -1580.4

$ws.Range("H33").Value = 554.0909
$ws.Range("I33").Value = 595.6111
$ws.Range("J33").Value = 367.25
$ws.Range("K33").Value = 595.6111
$ws.Range("L33").Value = 367.25
$ws.Range("M33").Value = -366.6111
$ws.Range("N33").Value = -825.25

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H40").Value = 1758.2113
$ws.Range("I40").Value = 1003.38983
$ws.Range("J40").Value = 5469.4165
$ws.Range("K40").Value = 1003.38983
$ws.Range("L40").Value = 5469.4165
$ws.Range("M40").Value = -828.38983
$ws.Range("N40").Value = -5819.4165

$ws.Range("H55").Value = 122.1875
$ws.Range("I55").Value = 138.15384
$ws.Range("J55").Value = 53
$ws.Range("K55").Value = 138.15384
$ws.Range("L55").Value = 53
$ws.Range("M55").Value = 75.84616
$ws.Range("N55").Value = -481

$ws.Range("H129").Value = 1042.6976
$ws.Range("J129").Value = 1151.579
$ws.Range("L129").Value = 3454.737
$ws.Range("N129").Value = -13454.737

$ws.Range("H131").Value = 2125
$ws.Range("I131").Value = 250
$ws.Range("K131").Value = 750
$ws.Range("M131").Value = 4290

$ws.Range("H132").Value = 2872.5881
$ws.Range("I132").Value = 1533.1923
$ws.Range("K132").Value = 4599.5769
$ws.Range("M132").Value = -2069.5769

$ws.Range("H137").Value = 3451.923
$ws.Range("I137").Value = 3380.4348
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 10141.3044
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -7591.304400000001
$ws.Range("N137").Value = -17100

$ws.Range("H138").Value = 2254.0667
$ws.Range("I138").Value = 1298.25
$ws.Range("J138").Value = 3346.4285
$ws.Range("K138").Value = 3894.75
$ws.Range("L138").Value = 10039.2855
$ws.Range("M138").Value = 1245.25
$ws.Range("N138").Value = -20319.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 33056.25
$ws.Range("J104").Value = 33056.25
$ws.Range("L104").Value = 33056.25
$ws.Range("N104").Value = -40044.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 602.375
$ws.Range("I64").Value = 573.3333
$ws.Range("J64").Value = 650.7778
$ws.Range("K64").Value = 573.3333
$ws.Range("L64").Value = 650.7778
$ws.Range("M64").Value = -348.3333
$ws.Range("N64").Value = -1100.7778

$ws.Range("H67").Value = 602.375
$ws.Range("I67").Value = 573.3333
$ws.Range("J67").Value = 650.7778
$ws.Range("K67").Value = 573.3333
$ws.Range("L67").Value = 650.7778
$ws.Range("M67").Value = 206.6667
$ws.Range("N67").Value = -2210.7778

$ws.Range("H86").Value = 3363.9285
$ws.Range("I86").Value = 4325.294
$ws.Range("J86").Value = 1878.1818
$ws.Range("K86").Value = 4325.294
$ws.Range("L86").Value = 1878.1818
$ws.Range("M86").Value = -3202.294
$ws.Range("N86").Value = -4124.1818

$ws.Range("H89").Value = 3363.9285
$ws.Range("I89").Value = 4325.294
$ws.Range("J89").Value = 1878.1818
$ws.Range("K89").Value = 21626.47
$ws.Range("L89").Value = 9390.909
$ws.Range("M89").Value = -16010.47
$ws.Range("N89").Value = -20622.909

$ws.Range("H99").Value = 50000
$ws.Range("I99").Value = 50000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 50000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -48502
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 234.27272
$ws.Range("I22").Value = 99.625
$ws.Range("J22").Value = 593.3333
$ws.Range("K22").Value = 99.625
$ws.Range("L22").Value = 593.3333
$ws.Range("M22").Value = 250.375
$ws.Range("N22").Value = -1293.3333

$ws.Range("H31").Value = 2360.3171
$ws.Range("I31").Value = 1258.0294
$ws.Range("K31").Value = 1258.0294
$ws.Range("M31").Value = -963.0293999999999

$ws.Range("H34").Value = 2360.3171
$ws.Range("I34").Value = 1258.0294
$ws.Range("K34").Value = 1258.0294
$ws.Range("M34").Value = -1056.0294

$ws.Range("H99").Value = 502228
$ws.Range("I99").Value = 1000012
$ws.Range("J99").Value = 4444
$ws.Range("K99").Value = 1000012
$ws.Range("L99").Value = 4444
$ws.Range("M99").Value = -998514
$ws.Range("N99").Value = -7440

$ws.Range("H126").Value = 502228
$ws.Range("I126").Value = 1000012
$ws.Range("J126").Value = 4444
$ws.Range("K126").Value = 3000036
$ws.Range("L126").Value = 13332
$ws.Range("M126").Value = -2997566
$ws.Range("N126").Value = -18272

$ws.Range("H134").Value = 1854.52
$ws.Range("I134").Value = 1069.5
$ws.Range("J134").Value = 2579.1538
$ws.Range("K134").Value = 3208.5
$ws.Range("L134").Value = 7737.4614
$ws.Range("M134").Value = -673.5
$ws.Range("N134").Value = -12807.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 43.57143
$ws.Range("I12").Value = 26
$ws.Range("J12").Value = 48.363636
$ws.Range("K12").Value = 78
$ws.Range("L12").Value = 145.090908
$ws.Range("M12").Value = 95
$ws.Range("N12").Value = -491.090908

$ws.Range("H38").Value = 42.636364
$ws.Range("I38").Value = 29.333334
$ws.Range("K38").Value = 88.00000199999999
$ws.Range("M38").Value = 258.999998

$ws.Range("H40").Value = 146.85
$ws.Range("I40").Value = 75.8
$ws.Range("J40").Value = 360
$ws.Range("K40").Value = 303.2
$ws.Range("L40").Value = 1440
$ws.Range("M40").Value = -234.2
$ws.Range("N40").Value = -1578

$ws.Range("H92").Value = 62500336
$ws.Range("I92").Value = 71428820
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 214286460
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -214285212
$ws.Range("N92").Value = -5496

$ws.Range("H97").Value = 268
$ws.Range("I97").Value = 146.16667
$ws.Range("J97").Value = 372.42856
$ws.Range("K97").Value = 438.50001
$ws.Range("L97").Value = 1117.28568
$ws.Range("M97").Value = 57.49998999999997
$ws.Range("N97").Value = -2109.28568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 29500
$ws.Range("J68").Value = 29500
$ws.Range("L68").Value = 29500
$ws.Range("N68").Value = -31122

$ws.Range("H71").Value = 29500
$ws.Range("J71").Value = 29500
$ws.Range("L71").Value = 88500
$ws.Range("N71").Value = -96612

$ws.Range("H102").Value = 6145.6665
$ws.Range("I102").Value = 3762
$ws.Range("J102").Value = 7337.5
$ws.Range("K102").Value = 3762
$ws.Range("L102").Value = 7337.5
$ws.Range("M102").Value = -2140
$ws.Range("N102").Value = -10581.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 769.95
$ws.Range("I22").Value = 559.2222
$ws.Range("J22").Value = 942.36365
$ws.Range("K22").Value = 559.2222
$ws.Range("L22").Value = 942.36365
$ws.Range("M22").Value = -264.2222
$ws.Range("N22").Value = -1532.36365

$ws.Range("H27").Value = 769.95
$ws.Range("I27").Value = 559.2222
$ws.Range("J27").Value = 942.36365
$ws.Range("K27").Value = 559.2222
$ws.Range("L27").Value = 942.36365
$ws.Range("M27").Value = -452.2222
$ws.Range("N27").Value = -1156.36365

$ws.Range("H46").Value = 1311.875
$ws.Range("I46").Value = 1340
$ws.Range("K46").Value = 1340
$ws.Range("M46").Value = -1152

$ws.Range("H55").Value = 304.2143
$ws.Range("I55").Value = 344.5
$ws.Range("J55").Value = 274
$ws.Range("K55").Value = 344.5
$ws.Range("L55").Value = 274
$ws.Range("M55").Value = -171.5
$ws.Range("N55").Value = -620

$ws.Range("H122").Value = 2788.889
$ws.Range("I122").Value = 2200
$ws.Range("J122").Value = 2862.5
$ws.Range("K122").Value = 6600
$ws.Range("L122").Value = 8587.5
$ws.Range("M122").Value = -4150
$ws.Range("N122").Value = -13487.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 685.2857
$ws.Range("I107").Value = 566.3333
$ws.Range("J107").Value = 774.5
$ws.Range("K107").Value = 1698.9999
$ws.Range("L107").Value = 2323.5
$ws.Range("M107").Value = 221.0001
$ws.Range("N107").Value = -6163.5

$ws.Range("H113").Value = 254.09091
$ws.Range("I113").Value = 157.75
$ws.Range("J113").Value = 511
$ws.Range("K113").Value = 473.25
$ws.Range("L113").Value = 1533
$ws.Range("M113").Value = 1696.75
$ws.Range("N113").Value = -5873

$ws.Range("H122").Value = 3300.8
$ws.Range("I122").Value = 3501.3333
$ws.Range("K122").Value = 10503.9999
$ws.Range("M122").Value = -8053.999899999999
